# Daily attendance processing - clear the "Recorded By" admin name from
# column G for every data row, and shrink column G now that it no longer
# needs to fit "Administrator, Miss Dina Nasr".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers (data rows) whose "Recorded By" (column G) cell currently
# holds "Administrator, Miss Dina Nasr" and must be cleared.
$rows = @(
    2,3,4,5,8,9,
    30,31,32,33,36,37,
    58,59,60,61,64,65,
    86,87,88,89,92,93,
    114,115,116,117,120,121,127,
    142,143,144,145,148,149,155,
    170,171,172,173,176,177,
    198,199,200,201,204,205,
    226,227,228,229,232,233,
    254,255,256,257,260,261,
    282,283,284,285,288,289,
    310,311,312,313,316,317,
    338,339,340,341,344,345,
    366,367,368,369,372,373,
    394,395,396,397,400,401,
    422,423,424,425,428,429,
    450,451,452,453,456,457,
    478,479,480,481,484,485,
    506,507,508,509,512,513,
    534,535,536,537,540,541,
    562,563,564,565,568,569,
    590,591,592,593,596,597,
    618,619,620,621,624,625,
    646,647,648,649,652,653
)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).ClearContents()
}

# Column G ("Recorded By") no longer needs to be wide enough for the
# administrator's name - shrink it from width 31 down to 13.
$ws.Columns.Item(7).ColumnWidth = 12.17
